# Update "paises" (Pais) worksheet with refreshed COVID-19 stats and
# re-sorted country rows, per commit "Update countries & provincias Spain".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 03:35"

# --- Row 51 / 52: Serbia and Corea del Sur swap position + new figures -
$ws.Range("A51").Value = "Corea del Sur"
$ws.Range("B51").Value = 11344
$ws.Range("C51").Value = 79
$ws.Range("D51").Value = 10340
$ws.Range("E51").Value = 735
$ws.Range("H51").Value = 269

$ws.Range("A52").Value = "Serbia"
$ws.Range("B52").Value = 11275
$ws.Range("D52").Value = 6277
$ws.Range("E52").Value = 4758
$ws.Range("H52").Value = 240

# --- Row 72: Sudan updated figures -------------------------------------
$ws.Range("B72").Value = 4346
$ws.Range("C72").Value = 200
$ws.Range("D72").Value = 749
$ws.Range("E72").Value = 3402
$ws.Range("G72").Value = 11
$ws.Range("H72").Value = 195

# --- Row 73 / 74: Luxemburgo and Guatemala swap position + new figures -
$ws.Range("A73").Value = "Guatemala"
$ws.Range("B73").Value = 4145
$ws.Range("C73").Value = 191
$ws.Range("D73").Value = 493
$ws.Range("E73").Value = 3584
$ws.Range("G73").Value = 5
$ws.Range("H73").Value = 68

$ws.Range("A74").Value = "Luxemburgo"
$ws.Range("B74").Value = 4001
$ws.Range("D74").Value = 3791
$ws.Range("E74").Value = 100
$ws.Range("H74").Value = 110

# --- Row 89: Croacia (unchanged name) updated figures ------------------
$ws.Range("D89").Value = 891
$ws.Range("E89").Value = 1179
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = 39

# --- Row 99: Eslovaquia (unchanged name) updated figures ---------------
$ws.Range("D99").Value = 1474
$ws.Range("E99").Value = 8
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 22
